# Adds a new "params" worksheet at the end of the workbook, populates it with
# the measurement/condition summary table, and restores the original active
# sheet/selection so only the intended structural change is introduced.

$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet

# Add the new sheet after the last existing sheet so it lands at the end (tab
# order: sgs, jsmatrix, params) and gets the next sheetId (3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "params"

# --- Header row -----------------------------------------------------------
# Write the condition/summary headers (columns B..BC) first, then the index
# header "meting" in column A last, matching the shared-string insertion
# order of the source workbook.
$headerTail = @('n','prev','prev_se','min_size','condition11','condition12','condition13','condition21','condition22','condition23','condition31','condition32','condition33','condition41','condition42','condition43','condition51','condition52','condition53','condition61','condition62','condition63','condition71','condition72','condition81','condition82','condition83','condition91','condition92','condition101','condition102','condition103','condition111','condition112','condition113','condition121','condition122','condition131','condition132','condition133','condition141','condition151','condition152','condition161','condition162','condition171','condition172','condition173','condition181','condition182','condition191','condition192','condition201','condition202')
for ($i = 0; $i -lt $headerTail.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headerTail[$i]
}
$ws.Range("A1").Value = "meting"

# --- Data rows --------------------------------------------------------------
$data = @(
    @(2003,4402,0.5533848250795094,0.007493844485279851,220.1,0.2921646746347942,0.2632575757575757,0.2534653465346535,0.2921646746347942,0.2632575757575757,0.2806451612903226,0.2921646746347942,0.2805555555555556,0.2404092071611253,0.2921646746347942,0.2805555555555556,0.2724137931034483,0.3817204301075269,0.361800346220427,0.3455637091394263,0.3817204301075269,0.361800346220427,0.3513011152416357,0.7154178674351584,0.7984361424847958,0.3817204301075269,0.361800346220427,0.3586367157242448,0.7154178674351584,0.800952380952381,0.3817204301075269,0.361800346220427,0.3605324074074074,0.3817204301075269,0.3626716604244694,0.360479797979798,0.7154178674351584,0.7549103330486764,0.3817204301075269,0.361800346220427,0.3521524347212421,0.7154178674351584,0.3817204301075269,0.361800346220427,0.7154178674351584,0.7346153846153847,0.3817204301075269,0.3626716604244694,0.3589093214965124,0.7154178674351584,0.7360126083530338,0.4787657597876576,0.4510948905109489,0.7154178674351584,0.710955710955711),
    @(2005,3522,0.5019875070982396,0.008426234876814374,176.1,0.2480376766091052,0.2310924369747899,0.2195652173913044,0.2480376766091052,0.2310924369747899,0.2024291497975708,0.2480376766091052,0.2361563517915309,0.212962962962963,0.2480376766091052,0.2361563517915309,0.208695652173913,0.3089214380825566,0.2896995708154507,0.2837729816147082,0.3089214380825566,0.2896995708154507,0.2887060583395662,0.7204116638078902,0.8132780082987552,0.3089214380825566,0.2896995708154507,0.2744186046511628,0.7204116638078902,0.8162583518930958,0.3089214380825566,0.2896995708154507,0.2914244186046512,0.3089214380825566,0.3014925373134328,0.299324831207802,0.7204116638078902,0.7753479125248509,0.3089214380825566,0.2896995708154507,0.2675276752767528,0.7204116638078902,0.3089214380825566,0.2896995708154507,0.7204116638078902,0.7448979591836735,0.3089214380825566,0.3014925373134328,0.2995495495495495,0.7204116638078902,0.7571569595261599,0.3938879456706282,0.3635509628933772,0.7204116638078902,0.7172727272727273),
    @(2007,4232,0.4468336483931947,0.007643271727985883,211.6,0.1554828150572831,0.1411483253588517,0.1306532663316583,0.1554828150572831,0.1411483253588517,0.148936170212766,0.1554828150572831,0.1393728222996516,0.11,0.1554828150572831,0.1393728222996516,0.1371681415929203,0.2216343327454439,0.20125,0.1906474820143885,0.2216343327454439,0.20125,0.1944818304172275,0.697560975609756,0.7483870967741936,0.2216343327454439,0.20125,0.196031746031746,0.697560975609756,0.7489139878366637,0.2216343327454439,0.20125,0.2020138451856514,0.2216343327454439,0.2080217539089055,0.2081911262798635,0.697560975609756,0.7364470391993327,0.2216343327454439,0.20125,0.1907514450867052,0.697560975609756,0.2216343327454439,0.20125,0.697560975609756,0.7037582903463523,0.2216343327454439,0.2080217539089055,0.2077922077922078,0.697560975609756,0.7044072948328267,0.3181980693600286,0.2906116088819634,0.697560975609756,0.6946508172362555),
    @(2009,3560,0.377808988764045,0.008127070281522939,178,0.1053484602917342,0.08669354838709678,0.0778688524590164,0.1053484602917342,0.08669354838709678,0.09642857142857143,0.1053484602917342,0.09437086092715231,0.07255520504731862,0.1053484602917342,0.09437086092715231,0.06389776357827476,0.1603206412825651,0.1433591004919185,0.133384734001542,0.1603206412825651,0.1433591004919185,0.1386430678466077,0.638228055783429,0.6972034715525555,0.1603206412825651,0.1433591004919185,0.1363265306122449,0.638228055783429,0.698477157360406,0.1603206412825651,0.1433591004919185,0.144793152639087,0.1603206412825651,0.1482300884955752,0.1477104874446086,0.638228055783429,0.6827852998065764,0.1603206412825651,0.1433591004919185,0.1376481312670921,0.638228055783429,0.1603206412825651,0.1433591004919185,0.638228055783429,0.6506968641114983,0.1603206412825651,0.1482300884955752,0.147819660014782,0.638228055783429,0.6575091575091575,0.2422041862451944,0.2161290322580645,0.638228055783429,0.636675235646958),
    @(2011,4375,0.3494857142857143,0.007209476721260907,218.75,0.08011049723756906,0.07524752475247524,0.07489878542510121,0.08011049723756906,0.07524752475247524,0.06853582554517133,0.08011049723756906,0.07845934379457917,0.06417112299465241,0.08011049723756906,0.07845934379457917,0.04984423676012461,0.1420794774088187,0.1314580941446613,0.1213333333333333,0.1420794774088187,0.1314580941446613,0.1280148423005566,0.6369426751592356,0.6826051112943117,0.1420794774088187,0.1314580941446613,0.1312217194570136,0.6369426751592356,0.6885688568856886,0.1420794774088187,0.1314580941446613,0.1310507674144038,0.1420794774088187,0.1292170591979631,0.1289087428206765,0.6369426751592356,0.6565377532228361,0.1420794774088187,0.1314580941446613,0.1292824822236587,0.6369426751592356,0.1420794774088187,0.1314580941446613,0.6369426751592356,0.6434456928838951,0.1420794774088187,0.1292170591979631,0.129156010230179,0.6369426751592356,0.6481774960380349,0.212356515867657,0.1887477313974592,0.6369426751592356,0.6384439359267735),
    @(2013,3341,0.2648907512720742,0.007635475844113779,167.05,0.03839732888146911,0.03023255813953488,0.02857142857142857,0.03839732888146911,0.03023255813953488,0.01915708812260536,0.03839732888146911,0.03747870528109029,0.02160493827160494,0.03839732888146911,0.03747870528109029,0.0316622691292876,0.07313540912382331,0.06566037735849056,0.05263157894736842,0.07313540912382331,0.06566037735849056,0.05709624796084829,0.5306691449814126,0.5651697699890471,0.07313540912382331,0.06566037735849056,0.05421103581800581,0.5306691449814126,0.5649122807017544,0.07313540912382331,0.06566037735849056,0.065284178187404,0.07313540912382331,0.05756013745704467,0.05598621877691645,0.5306691449814126,0.5506257110352674,0.07313540912382331,0.06566037735849056,0.06143344709897611,0.5306691449814126,0.07313540912382331,0.06566037735849056,0.5306691449814126,0.5401234567901234,0.07313540912382331,0.05756013745704467,0.05613126079447323,0.5306691449814126,0.537180910099889,0.1386313465783665,0.1164772727272727,0.5306691449814126,0.5284872298624754),
    @(2015,4322,0.2320684868116613,0.006422101970115271,216.1,0.0467065868263473,0.03565640194489465,0.02666666666666667,0.0467065868263473,0.03565640194489465,0.0325,0.0467065868263473,0.03593556381660471,0.01590909090909091,0.0467065868263473,0.03593556381660471,0.02133333333333333,0.0733652312599681,0.06030150753768844,0.05276381909547739,0.0733652312599681,0.06030150753768844,0.05562060889929742,0.4705014749262537,0.5190972222222222,0.0733652312599681,0.06030150753768844,0.05145565335138795,0.4705014749262537,0.5285977859778598,0.0733652312599681,0.06030150753768844,0.05938375350140056,0.0733652312599681,0.06430288461538461,0.06332931242460796,0.4705014749262537,0.5138023152270703,0.0733652312599681,0.06030150753768844,0.05520304568527919,0.4705014749262537,0.0733652312599681,0.06030150753768844,0.4705014749262537,0.4748822605965463,0.0733652312599681,0.06430288461538461,0.06238643246517262,0.4705014749262537,0.4773289365210223,0.1230613621038436,0.1062455132806892,0.4705014749262537,0.4759316770186335),
    @(2017,4145,0.2130277442702051,0.006360451156136092,207.25,0.03243847874720358,0.0202808112324493,0.01628664495114007,0.03243847874720358,0.0202808112324493,0.0131578947368421,0.03243847874720358,0.02686915887850467,0.02040816326530612,0.03243847874720358,0.02686915887850467,0.01470588235294118,0.0655226209048362,0.05774569683509161,0.05212968849332485,0.0655226209048362,0.05774569683509161,0.05322294500295683,0.4353954581049335,0.4878993223620523,0.0655226209048362,0.05774569683509161,0.0499001996007984,0.4353954581049335,0.4948559670781893,0.0655226209048362,0.05774569683509161,0.05752961082910321,0.0655226209048362,0.05765765765765766,0.05693519079345851,0.4353954581049335,0.4768339768339768,0.0655226209048362,0.05774569683509161,0.05800604229607251,0.4353954581049335,0.0655226209048362,0.05774569683509161,0.4353954581049335,0.4428807947019868,0.0655226209048362,0.05765765765765766,0.05643203883495146,0.4353954581049335,0.4508414526129318,0.1140167364016736,0.09418070444104135,0.4353954581049335,0.4383223684210527),
    @(2019,3486,0.2240390131956397,0.007062856420827443,174.3,0.05236270753512133,0.04528985507246377,0.03802281368821293,0.05236270753512133,0.04528985507246377,0.0325,0.05236270753512133,0.04582210242587601,0.03551912568306011,0.05236270753512133,0.04582210242587601,0.03865979381443299,0.0762660158633313,0.06723237597911227,0.06113207547169811,0.0762660158633313,0.06723237597911227,0.06388888888888888,0.4488107549120993,0.5050890585241731,0.0762660158633313,0.06723237597911227,0.06354249404289118,0.4488107549120993,0.5213793103448275,0.0762660158633313,0.06723237597911227,0.06737120211360634,0.0762660158633313,0.06874557051736357,0.06743185078909612,0.4488107549120993,0.4889466840052016,0.0762660158633313,0.06723237597911227,0.06423982869379015,0.4488107549120993,0.0762660158633313,0.06723237597911227,0.4488107549120993,0.4539400665926748,0.0762660158633313,0.06874557051736357,0.06647398843930635,0.4488107549120993,0.4652532391048292,0.1377530766177054,0.1181619256017506,0.4488107549120993,0.4580573951434879)
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# --- Formatting -------------------------------------------------------------
# Column A ("meting") and the header row both use the bold/centered/bordered
# header style already defined in the workbook (style index 1, as seen on the
# "sgs" and "jsmatrix" sheets). Copy it over instead of re-declaring bold /
# border / alignment so no new style entries are created.
$styleSource = $wb.Worksheets.Item(1).Range("B1")
$styleSource.Copy()
$ws.Range("A1:BC1").PasteSpecial(-4122)
$ws.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the original active sheet/selection so "params" isn't left as the
# selected tab.
$originalActive.Select()
